# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures (Universalis refresh) to the
# per-job Leve profit tables (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR), matching the
# scheduled-runner data refresh described in the commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 3571621
$ws.Range("I6").Value = 4081709.8
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 12245129.4
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = -12245017.4
$ws.Range("N6").Value = -3224

$ws.Range("H15").Value = 1333387.8
$ws.Range("I15").Value = 1333387.8
$ws.Range("K15").Value = 4000163.4
$ws.Range("M15").Value = -3999994.4

$ws.Range("H92").Value = 370
$ws.Range("I92").Value = 370
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 370
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 878
$ws.Range("N92").ClearContents()

$ws.Range("H96").Value = 274.7647
$ws.Range("I96").Value = 226.3077
$ws.Range("J96").Value = 432.25
$ws.Range("K96").Value = 678.9231
$ws.Range("L96").Value = 1296.75
$ws.Range("M96").Value = 694.0769
$ws.Range("N96").Value = -4042.75

$ws.Range("H98").Value = 2408
$ws.Range("I98").Value = 2521.1765
$ws.Range("J98").Value = 1766.6666
$ws.Range("K98").Value = 2521.1765
$ws.Range("L98").Value = 1766.6666
$ws.Range("M98").Value = -1023.1765
$ws.Range("N98").Value = -4762.6666

$ws.Range("H100").Value = 1876.9286
$ws.Range("I100").Value = 1839.5714
$ws.Range("K100").Value = 1839.5714
$ws.Range("M100").Value = -1298.5714

$ws.Range("H116").Value = 1686493.4
$ws.Range("I116").Value = 3000.6667
$ws.Range("J116").Value = 3706684.5
$ws.Range("K116").Value = 3000.6667
$ws.Range("L116").Value = 3706684.5
$ws.Range("M116").Value = 441.3332999999998
$ws.Range("N116").Value = -3713568.5

$ws.Range("H122").Value = 2408
$ws.Range("I122").Value = 2521.1765
$ws.Range("J122").Value = 1766.6666
$ws.Range("K122").Value = 7563.529500000001
$ws.Range("L122").Value = 5299.9998
$ws.Range("M122").Value = -5113.529500000001
$ws.Range("N122").Value = -10199.9998

$ws.Range("H138").Value = 2625.8538
$ws.Range("I138").Value = 1831.1072
$ws.Range("J138").Value = 3037.9443
$ws.Range("K138").Value = 5493.321599999999
$ws.Range("L138").Value = 9113.832900000001
$ws.Range("M138").Value = -353.3215999999993
$ws.Range("N138").Value = -19393.8329

$ws.Range("H141").Value = 6739.7915
$ws.Range("I141").Value = 1676.3158
$ws.Range("J141").Value = 25981
$ws.Range("K141").Value = 5028.9474
$ws.Range("L141").Value = 77943
$ws.Range("M141").Value = 151.0526
$ws.Range("N141").Value = -88303


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1600
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 1600
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 1600
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -2412

$ws.Range("H91").Value = 1600
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 1600
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 1600
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -4408

$ws.Range("H102").Value = 2828
$ws.Range("I102").Value = 1704
$ws.Range("J102").Value = 3952
$ws.Range("K102").Value = 1704
$ws.Range("L102").Value = 3952
$ws.Range("M102").Value = -82
$ws.Range("N102").Value = -7196


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2541.5
$ws.Range("I86").Value = 4023.5715
$ws.Range("J86").Value = 1598.3636
$ws.Range("K86").Value = 4023.5715
$ws.Range("L86").Value = 1598.3636
$ws.Range("M86").Value = -2900.5715
$ws.Range("N86").Value = -3844.3636

$ws.Range("H89").Value = 2541.5
$ws.Range("I89").Value = 4023.5715
$ws.Range("J89").Value = 1598.3636
$ws.Range("K89").Value = 20117.8575
$ws.Range("L89").Value = 7991.817999999999
$ws.Range("M89").Value = -14501.8575
$ws.Range("N89").Value = -19223.818

$ws.Range("H94").Value = 2745.4546
$ws.Range("I94").Value = 2825
$ws.Range("J94").Value = 2533.3333
$ws.Range("K94").Value = 2825
$ws.Range("L94").Value = 2533.3333
$ws.Range("M94").Value = -2374
$ws.Range("N94").Value = -3435.3333

$ws.Range("H99").Value = 1771.7646
$ws.Range("I99").Value = 1970
$ws.Range("J99").Value = 1710.7693
$ws.Range("K99").Value = 1970
$ws.Range("L99").Value = 1710.7693
$ws.Range("M99").Value = -472
$ws.Range("N99").Value = -4706.7693

$ws.Range("H104").Value = 34000
$ws.Range("J104").Value = 34000
$ws.Range("L104").Value = 34000
$ws.Range("N104").Value = -40988

$ws.Range("H105").Value = 2996.6667
$ws.Range("I105").Value = 2733.3333
$ws.Range("J105").Value = 3260
$ws.Range("K105").Value = 2733.3333
$ws.Range("L105").Value = 3260
$ws.Range("M105").Value = -986.3332999999998
$ws.Range("N105").Value = -6754


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2476.3196
$ws.Range("I31").Value = 1343.7885
$ws.Range("J31").Value = 3785.0222
$ws.Range("K31").Value = 1343.7885
$ws.Range("L31").Value = 3785.0222
$ws.Range("M31").Value = -1048.7885
$ws.Range("N31").Value = -4375.022199999999

$ws.Range("H34").Value = 2476.3196
$ws.Range("I34").Value = 1343.7885
$ws.Range("J34").Value = 3785.0222
$ws.Range("K34").Value = 1343.7885
$ws.Range("L34").Value = 3785.0222
$ws.Range("M34").Value = -1141.7885
$ws.Range("N34").Value = -4189.022199999999

$ws.Range("H62").Value = 9769
$ws.Range("I62").Value = 2871
$ws.Range("J62").Value = 13601.223
$ws.Range("K62").Value = 2871
$ws.Range("L62").Value = 13601.223
$ws.Range("M62").Value = -2247
$ws.Range("N62").Value = -14849.223

$ws.Range("H65").Value = 9769
$ws.Range("I65").Value = 2871
$ws.Range("J65").Value = 13601.223
$ws.Range("K65").Value = 14355
$ws.Range("L65").Value = 68006.11500000001
$ws.Range("M65").Value = -11235
$ws.Range("N65").Value = -74246.11500000001

$ws.Range("H132").Value = 2026.0476
$ws.Range("I132").Value = 1498
$ws.Range("J132").Value = 2606.9
$ws.Range("K132").Value = 4494
$ws.Range("L132").Value = 7820.700000000001
$ws.Range("M132").Value = -1964
$ws.Range("N132").Value = -12880.7


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 2600
$ws.Range("I97").Value = 5000
$ws.Range("J97").Value = 1400
$ws.Range("K97").Value = 15000
$ws.Range("L97").Value = 4200
$ws.Range("M97").Value = -14504
$ws.Range("N97").Value = -5192


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1859.6666
$ws.Range("I80").Value = 1862.375
$ws.Range("J80").Value = 1857.5
$ws.Range("K80").Value = 1862.375
$ws.Range("L80").Value = 1857.5
$ws.Range("M80").Value = -864.375
$ws.Range("N80").Value = -3853.5

$ws.Range("H83").Value = 1859.6666
$ws.Range("I83").Value = 1862.375
$ws.Range("J83").Value = 1857.5
$ws.Range("K83").Value = 9311.875
$ws.Range("L83").Value = 9287.5
$ws.Range("M83").Value = -4319.875
$ws.Range("N83").Value = -19271.5

$ws.Range("H97").Value = 820
$ws.Range("I97").Value = 821.4286
$ws.Range("J97").Value = 800
$ws.Range("K97").Value = 821.4286
$ws.Range("L97").Value = 800
$ws.Range("M97").Value = -325.4286
$ws.Range("N97").Value = -1792

$ws.Range("H113").Value = 2280.182
$ws.Range("I113").Value = 2063
$ws.Range("K113").Value = 2063
$ws.Range("M113").Value = 107


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4166.5557
$ws.Range("I7").Value = 4166.5
$ws.Range("J7").Value = 4166.6665
$ws.Range("K7").Value = 4166.5
$ws.Range("L7").Value = 4166.6665
$ws.Range("M7").Value = -4054.5
$ws.Range("N7").Value = -4390.6665

$ws.Range("H40").Value = 33336420
$ws.Range("I40").Value = 58826076
$ws.Range("J40").Value = 3793.7693
$ws.Range("K40").Value = 58826076
$ws.Range("L40").Value = 3793.7693
$ws.Range("M40").Value = -58825940
$ws.Range("N40").Value = -4065.7693

$ws.Range("H68").Value = 2292.8
$ws.Range("I68").Value = 1999.2727
$ws.Range("J68").Value = 3100
$ws.Range("K68").Value = 1999.2727
$ws.Range("L68").Value = 3100
$ws.Range("M68").Value = -1250.2727
$ws.Range("N68").Value = -4598

$ws.Range("H71").Value = 2292.8
$ws.Range("I71").Value = 1999.2727
$ws.Range("J71").Value = 3100
$ws.Range("K71").Value = 9996.363499999999
$ws.Range("L71").Value = 15500
$ws.Range("M71").Value = -6252.363499999999
$ws.Range("N71").Value = -22988

$ws.Range("H97").Value = 17523.182
$ws.Range("J97").Value = 17523.182
$ws.Range("L97").Value = 17523.182
$ws.Range("N97").Value = -19505.182

$ws.Range("H100").Value = 3319.8
$ws.Range("I100").Value = 2942.875
$ws.Range("J100").Value = 3571.0833
$ws.Range("K100").Value = 2942.875
$ws.Range("L100").Value = 3571.0833
$ws.Range("M100").Value = -2401.875
$ws.Range("N100").Value = -4653.0833

$ws.Range("H126").Value = 4166.5557
$ws.Range("I126").Value = 4166.5
$ws.Range("J126").Value = 4166.6665
$ws.Range("K126").Value = 12499.5
$ws.Range("L126").Value = 12499.9995
$ws.Range("M126").Value = -10029.5
$ws.Range("N126").Value = -17439.9995


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4370
$ws.Range("I62").Value = 3327.3333
$ws.Range("J62").Value = 5934
$ws.Range("K62").Value = 3327.3333
$ws.Range("L62").Value = 5934
$ws.Range("M62").Value = -2703.3333
$ws.Range("N62").Value = -7182

$ws.Range("H65").Value = 4370
$ws.Range("I65").Value = 3327.3333
$ws.Range("J65").Value = 5934
$ws.Range("K65").Value = 16636.6665
$ws.Range("L65").Value = 29670
$ws.Range("M65").Value = -13516.6665
$ws.Range("N65").Value = -35910

$ws.Range("H81").Value = 10375.125
$ws.Range("I81").Value = 8500.166999999999
$ws.Range("K81").Value = 17000.334
$ws.Range("M81").Value = -15939.334

$ws.Range("H84").Value = 10375.125
$ws.Range("I84").Value = 8500.166999999999
$ws.Range("K84").Value = 85001.67
$ws.Range("M84").Value = -79697.67

$ws.Range("H96").Value = 3276
$ws.Range("I96").Value = 2800
$ws.Range("J96").Value = 3344
$ws.Range("K96").Value = 2800
$ws.Range("L96").Value = 3344
$ws.Range("M96").Value = -1427
$ws.Range("N96").Value = -6090

$ws.Range("H122").Value = 2361.5417
$ws.Range("I122").Value = 2037.2667
$ws.Range("J122").Value = 2902
$ws.Range("K122").Value = 6111.800099999999
$ws.Range("L122").Value = 8706
$ws.Range("M122").Value = -3661.800099999999
$ws.Range("N122").Value = -13606

$ws.Range("H136").Value = 1848.5088
$ws.Range("I136").Value = 1845.0238
$ws.Range("J136").Value = 1858.2667
$ws.Range("K136").Value = 5535.0714
$ws.Range("L136").Value = 5574.800099999999
$ws.Range("M136").Value = -2985.0714
$ws.Range("N136").Value = -10674.8001
